$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting the existing B (query) and
# C (dbExcel filename) columns one place to the right. Excel's own
# "Insert" behaviour copies the cell formatting (including the wrap-text
# style used in row 2) from the column to the left, which is what the
# target workbook shows for the new column.
$ws.Columns("B:B").Insert()

# Match column A's width on the newly inserted column as closely as this
# runtime's rounding allows (Excel normally clones the left neighbour's
# width automatically on insert).
$ws.Columns("B:B").ColumnWidth = 75

# Populate the new "StatQuery" column.
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''Saint Bernard'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
